$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (CT1) - fill in execution data (order: TRACEID, WORKFLOWINSTANCEID, TIMESTAMP)
$ws.Range("H2").Value = "a2735f97b3d3ef68"
$ws.Range("I2").Value = "2.16.840.1.113883.2.9.2.150.4.4.76b2b73e459caf4a8471d4f4179edfff8ba13c710e1f16c0cfa9557c637996eb.be36b794d5^^^^urn:ihe:iti:xdw:2013:workflowInstanceId"
$ws.Range("G2").Value = "2023-02-03T11:47:29Z"
$ws.Range("F2").Value = 44960
$ws.Range("K2").Value = "SI"
$ws.Range("L2").Value = "PASS"

# Row 3 (CT2)
$ws.Range("H3").Value = "32bda0909e8f1ca5"
$ws.Range("I3").Value = "2.16.840.1.113883.2.9.2.150.4.4.b6fe5ef53219755dd2bc9278b25283017acb89fb965a5b3dcb9ce0e460b832e7.41998613c8^^^^urn:ihe:iti:xdw:2013:workflowInstanceId"
$ws.Range("G3").Value = "2023-02-03T12:16:55Z"
$ws.Range("F3").Value = 44960
$ws.Range("K3").Value = "SI"
$ws.Range("L3").Value = "PASS"

# Row 4 (CT3)
$ws.Range("H4").Value = "bd7534cc9da542e3"
$ws.Range("I4").Value = "2.16.840.1.113883.2.9.2.150.4.4.dc57145a2f00d2bbdc51b87566ab588a0369dfa4b95e4ceb4dc8921eba005d7f.ab75a232d7^^^^urn:ihe:iti:xdw:2013:workflowInstanceId"
$ws.Range("G4").Value = "2023-02-03T12:26:01Z"
$ws.Range("F4").Value = 44960
$ws.Range("K4").Value = "SI"
$ws.Range("L4").Value = "PASS"

# Row 5 (CT4)
$ws.Range("H5").Value = "9bc547a176d7874c"
$ws.Range("I5").Value = "2.16.840.1.113883.2.9.2.150.4.4.00c8b21447e5435c0879e69ba2e30aa77b77940f9f7cb24bd89f79cc005d7db9.9ff9e5fdb7^^^^urn:ihe:iti:xdw:2013:workflowInstanceId"
$ws.Range("G5").Value = "2023-02-03T12:31:26Z"
$ws.Range("F5").Value = 44960
$ws.Range("K5").Value = "SI"
$ws.Range("L5").Value = "PASS"

# Row 6 (CT5)
$ws.Range("H6").Value = "44cd7a6254883840"
$ws.Range("I6").Value = "2.16.840.1.113883.2.9.2.150.4.4.07b01a2a3a374c407cfc2f18cb9339b72575b8c0d933476fbed72d4a550b701e.8fa3938ccf^^^^urn:ihe:iti:xdw:2013:workflowInstanceId"
$ws.Range("G6").Value = "2023-02-03T12:36:22Z"
$ws.Range("F6").Value = 44960
$ws.Range("K6").Value = "SI"
$ws.Range("L6").Value = "PASS"

# Row 7 - esito K7 unchanged value "NO" (kept)
$ws.Range("K7").Value = "NO"

# Row 8 (VALIDAZIONE_TOKEN_JWT_CAMPO_LAB_KO) - fill in execution data
$ws.Range("H8").Value = "e85f99f5eca606ee"
$ws.Range("I8").Value = "UNKNOWN_WORKFLOW_ID"
$ws.Range("G8").Value = "2023-02-03T11:44:35Z"
$ws.Range("J8").Value = "Il referto non viene prodotto e viene visualizzato a video il seguente messaggio di errore: 'Errore durante la validazione del referto: Campo token JWT non valido.'"
$ws.Range("F8").Value = 44960
$ws.Range("K8").Value = "SI"
$ws.Range("N8").Value = "SI"

# Row 9 (VALIDAZIONE_LAB_TIMEOUT)
$ws.Range("K9").Value = "SI"
$ws.Range("N9").Value = "SI"

# Row 10 (CT6_KO)
$ws.Range("K10").Value = "NO"
$ws.Range("N10").Value = "SI"

# Row 11 (CT7_KO) - fill in execution data (I11 & J11 reuse already-introduced strings)
$ws.Range("H11").Value = "e038f83af4961d44"
$ws.Range("I11").Value = "UNKNOWN_WORKFLOW_ID"
$ws.Range("G11").Value = "2023-02-03T15:58:22Z"
$ws.Range("J11").Value = "Il referto non viene prodotto e viene visualizzato a video il seguente messaggio di errore: 'Errore durante la validazione del referto: Campo token JWT non valido.'"
$ws.Range("F11").Value = 44960
$ws.Range("K11").Value = "SI"
$ws.Range("N11").Value = "SI"

# Row 17 (CT13_KO) - fill in execution data
$ws.Range("H17").Value = "ab05e1350ac2cb64"
$ws.Range("I17").Value = "2.16.840.1.113883.2.9.2.150.4.4.6999035a578ad82dcb97ebaa144285e3771d4cd3aedad75ae60b623f1eaaff4c.3b5025ead1^^^^urn:ihe:iti:xdw:2013:workflowInstanceId"
$ws.Range("G17").Value = "2023-02-03T16:37:50Z"
$ws.Range("J17").Value = "Il referto non viene prodotto e viene visualizzato a video il seguente messaggio di errore: 'Errore durante la validazione del referto: Errore semantico.'"
$ws.Range("F17").Value = 44960
$ws.Range("K17").Value = "SI"
$ws.Range("N17").Value = "SI"

# Row 12 (CT8_KO) - clear previous execution data, add new note
$ws.Range("F12").ClearContents()
$ws.Range("H12").ClearContents()
$ws.Range("I12").ClearContents()
$ws.Range("M12").Value = "Il campo oggetto del caso di test è sempre valorizzato con 'N' o 'V'. Non sono possibili valori differenti."
$ws.Range("K12").Value = "NO"
$ws.Range("N12").Value = "SI"

# Row 18 (CT14_KO) - fill in execution data (J18 reuses already-introduced string)
$ws.Range("H18").Value = "e75474093f38aeed"
$ws.Range("I18").Value = "2.16.840.1.113883.2.9.2.150.4.4.6999035a578ad82dcb97ebaa144285e3771d4cd3aedad75ae60b623f1eaaff4c.c640424e29^^^^urn:ihe:iti:xdw:2013:workflowInstanceId"
$ws.Range("G18").Value = "2023-02-03T16:58:34Z"
$ws.Range("J18").Value = "Il referto non viene prodotto e viene visualizzato a video il seguente messaggio di errore: 'Errore durante la validazione del referto: Errore semantico.'"
$ws.Range("F18").Value = 44960
$ws.Range("K18").Value = "SI"
$ws.Range("N18").Value = "SI"

# Row 13 (CT9_KO)
$ws.Range("K13").Value = "NO"
$ws.Range("N13").Value = "SI"

# Row 14 (CT10_KO)
$ws.Range("K14").Value = "NO"
$ws.Range("N14").Value = "SI"

# Row 15 (CT11_KO)
$ws.Range("K15").Value = "NO"
$ws.Range("N15").Value = "SI"

# Row 16 (CT12_KO) - clear previous execution data, add new note (reuses existing note string)
$ws.Range("F16").ClearContents()
$ws.Range("H16").ClearContents()
$ws.Range("I16").ClearContents()
$ws.Range("M16").Value = "Il campo oggetto del caso di test è sempre valorizzato in maniera corretta."
$ws.Range("K16").Value = "NO"
$ws.Range("N16").Value = "SI"

# Row 19 (CT15_KO)
$ws.Range("K19").Value = "NO"
$ws.Range("N19").Value = "SI"

# Row 20 (CT16_KO)
$ws.Range("K20").Value = "NO"
$ws.Range("N20").Value = "SI"

# Set view: scroll so F10 is top-left, and select N15
$ws.Range("F10").Select()
$excel.ActiveWindow.ScrollRow = 10
$excel.ActiveWindow.ScrollColumn = 6
$ws.Range("N15").Select()
